# PSRrysowanieDokumentacja.docx edit script
# Reproduces:
#  1) "W następnej kolejności musimy włączyć moduł klienta."
#       -> "Następnie włączamy moduły klientów."        (split into several runs)
#  2) an extra formatted space inserted after "...poniżej."
#  3) a re-split of the "Nieco niżej / umiejscowione..." run pair at a new
#     text boundary (no net text change)
#  4) a re-split of " czy *.png. " around "png" (no net text change)

$d = $word.ActiveDocument

function Split-RunAt($absPos) {
    # Force a run boundary immediately before the character at $absPos by
    # toggling a character formatting attribute on the preceding character
    # and then reverting it. The engine keeps runs split at positions whose
    # formatting was explicitly (even if transiently) touched, even though
    # the final effective formatting is unchanged.
    $touch = $d.Range($absPos - 1, $absPos)
    $touch.Bold = 1
    $touch.Bold = 0
}

function Find-Bounds($needle, $searchStart, $searchEnd) {
    $scope = $d.Range($searchStart, $searchEnd)
    $ok = $scope.Find.Execute($needle)
    if (-not $ok) {
        throw "Not found: $needle"
    }
    return $scope
}

# ---------------------------------------------------------------------
# Change 1: replace the sentence, then split the replacement into runs
# ---------------------------------------------------------------------

$para = $d.Content
$para.Find.Execute("Przed uruchomieniem")
$paraStart = $para.Start

$ok = $d.Content.Find.Execute( `
    "W następnej kolejności musimy włączyć moduł klienta.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Następnie włączamy moduły klientów.", 2)
Write-Output "sentence replace ok=$ok"

# locate the freshly inserted sentence
$sentRange = $d.Content
$sentRange.Find.Execute("Następnie włączamy moduły klientów.")
$sentStart = $sentRange.Start
$sentEnd = $sentRange.End
Write-Output "sentence start=$sentStart end=$sentEnd"

# split boundaries, located by searching within the sentence range:
#   "Następnie" | " " | "włączamy" | " moduł" | "y" | " klient" | "ów" | ". Po..."
$bNastepnie = Find-Bounds "Następnie" $sentStart $sentEnd
Split-RunAt $bNastepnie.End

$bSpace1 = Find-Bounds " włączamy" $sentStart $sentEnd
Split-RunAt ($bSpace1.Start + 1)

$bWlaczamy = Find-Bounds "włączamy" $sentStart $sentEnd
Split-RunAt $bWlaczamy.End

$bModul = Find-Bounds " moduł" $sentStart $sentEnd
Split-RunAt $bModul.End

$bY = Find-Bounds "y klient" $sentStart $sentEnd
Split-RunAt ($bY.Start + 1)

$bKlient = Find-Bounds " klient" $sentStart $sentEnd
Split-RunAt $bKlient.End

$bOw = Find-Bounds "ów." $sentStart $sentEnd
Split-RunAt ($bOw.Start + 2)

# preserve the pre-existing run boundary right before " opisana została poniżej."
$bOpisana = Find-Bounds " opisana została poniżej." $sentStart ($sentEnd + 200)
Split-RunAt $bOpisana.Start

Write-Output "change 1 done"
